$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (old "1304060 - Maria das Graças de Almeida Felipe" row under
# "Docentes responsáveis:"), which shifts all subsequent rows up by one.
$ws.Rows.Item(13).Delete()

# After the shift, patch up the cells whose content was re-pointed at other
# (pre-existing) strings in the shared-strings table.
$ws.Range("B10").Value = "1304060 - Maria das Graças de Almeida Felipe"
$ws.Range("C10").Value = "1304060 - Maria das Graças de Almeida Felipe"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C15").Value = "01/01/2022"

$ws.Range("B18").Value = "1304060 - Maria das Graças de Almeida Felipe"
$ws.Range("C18").Value = "1304060 - Maria das Graças de Almeida Felipe"

$ws.Range("B19").Value = "Duas provas escritas (P1 e P2) distribuídas no semestre.Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0."
$ws.Range("C19").Value = "Duas provas escritas (P1 e P2) distribuídas no semestre.Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0."

$ws.Range("B20").Value = "MF = média finalMF = (P1 + P2)/2"
$ws.Range("C20").Value = "MF = média finalMF = (P1 + P2)/2"

$ws.Range("B21").Value = "Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0."
$ws.Range("C21").Value = "Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0."
